$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.808.55"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.892.63"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7950"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.89"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3164"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.44"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07044"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08076"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7685"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.895.19"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.345"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.49"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.836.94"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.992"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.86"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.51"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007700"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.357"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +20.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.0000"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.148.20"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1639"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.348"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.13"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.70"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.056"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.399"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.538"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.440"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05695"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.042"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.260"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7379"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9996"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.620"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01911"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.785"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4404"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.38"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.815"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8407"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.032.28"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.04"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.34%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.871"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.991"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.420"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.044.49"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.14%  "
